$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9999999936121029
$ws.Range("A2").Value = 0.99731590709603957
$ws.Range("A3").Value = 0.99408874165908656
$ws.Range("A4").Value = 0.99758702317637926
$ws.Range("A5").Value = 0.98795880092599475
$ws.Range("A6").Value = 0.96509516356814784
$ws.Range("A7").Value = 0.96256653162832984
$ws.Range("A8").Value = 0.95969733770111731
$ws.Range("A9").Value = 0.95866195589715952
$ws.Range("A10").Value = 0.95857665035849382
$ws.Range("A11").Value = 0.95868312688103341
$ws.Range("A12").Value = 0.95914248495169052
$ws.Range("A13").Value = 0.9633379582707895
$ws.Range("A14").Value = 0.96149483406614933
$ws.Range("A15").Value = 0.96136391695502987
$ws.Range("A16").Value = 0.96192896267180727
$ws.Range("A17").Value = 0.95822143513145464
$ws.Range("A18").Value = 0.95711254362214837
$ws.Range("A19").Value = 0.99305136606297817
$ws.Range("A20").Value = 0.97242920994331772
$ws.Range("A21").Value = 0.96695955291371327
$ws.Range("A22").Value = 0.95825644095146933
$ws.Range("A23").Value = 0.98705461664649452
$ws.Range("A24").Value = 0.97403442724328082
$ws.Range("A25").Value = 0.9675775473393613
$ws.Range("A26").Value = 0.9618984971960598
$ws.Range("A27").Value = 0.95705520377272402
$ws.Range("A28").Value = 0.93559306503399653
$ws.Range("A29").Value = 0.92032563798157008
$ws.Range("A30").Value = 0.91375605749217792
$ws.Range("A31").Value = 0.90610290441463026
$ws.Range("A32").Value = 0.90442361572421803
$ws.Range("A33").Value = 0.9039036142691248
